$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.358.57"
$ws.Range("E2").Value = "  +11.24%  "
$ws.Range("D3").Value = "1.877.46"
$ws.Range("E3").Value = "  +7.82%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9978"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "250.32"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9966"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +3.54%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "45.11"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +9.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2853"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +9.55%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06553"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +6.43%  "
$ws.Range("D11").Value = "1.871.67"
$ws.Range("E11").Value = "  +8.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "17.00"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +5.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07184"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.60%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6654"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +10.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "86.06"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +11.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.810"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +8.29%  "
$ws.Range("D17").Value = "30.378.58"
$ws.Range("E17").Value = "  +11.53%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9950"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.35%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007539"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +6.14%  "
$ws.Range("E20").Value = "  +9.93%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9972"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").Value = "2.096.99"
$ws.Range("E22").Value = "  +7.45%  "
$ws.Range("E23").Value = "  +6.66%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.506"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +7.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.994"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +7.11%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "144.31"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "135.47"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +26.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.73"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +9.68%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.937"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +6.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.400"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.59%  "
$ws.Range("E31").Value = "  +7.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08610"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +8.45%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.882"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +5.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05071"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.132"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +11.90%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6826"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +10.53%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.683"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +3.58%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.315"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +14.68%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.756"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9545"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01625"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.13%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.148"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +7.82%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.45"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +5.15%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9950"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4181"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +9.31%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.447"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1248"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +8.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05649"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +5.53%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.311"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "32.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +8.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.345"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +8.19%  "
